$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update batch number text
$ws.Range("C3").Value = "BATCH44444"

# Update quantities (request detail) for rows 11 and 12
$ws.Range("E11").Value = 5000
$ws.Range("E12").Value = 13000

# Move the active selection to C4 (matches the saved cursor position)
$ws.Range("C4").Select()
